$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data (row -> A, B label, C, D, E)
# Rows 2-9 are line1..line8, rows 10-17 are extr1..extr8
$rows = @(
    @{ Row = 2;  A = 0;  B = "line1"; C = 7;  D = 9;  E = $true  },
    @{ Row = 3;  A = 1;  B = "line2"; C = 9;  D = 8;  E = $false },
    @{ Row = 4;  A = 2;  B = "line3"; C = 8;  D = 10; E = $true  },
    @{ Row = 5;  A = 3;  B = "line4"; C = 8;  D = 11; E = $true  },
    @{ Row = 6;  A = 4;  B = "line5"; C = 10; D = 5;  E = $true  },
    @{ Row = 7;  A = 5;  B = "line6"; C = 12; D = 8;  E = $true  },
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $true  },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $true  },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $false },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $true  },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $true  }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
}

# Rows 16 and 17 are new - copy formatting (the "s=1" style) for column A from an existing
# styled cell (A2) so the new rows match the rest of the column A formatting.
$ws.Range("A2").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
